# Update countries & provincias Spain
# -----------------------------------
# Refresh the COVID "Pais" table with the 15:52 snapshot (previous data was
# timestamped 15:22). Most countries keep their row position and simply get
# refreshed totals; a handful of countries (Singapur, Islandia, Eslovaquia)
# overtook their neighbours in "Casos totales" and so their rows now carry
# the updated figures while the countries they passed slide down to absorb
# the values that used to sit above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 10 de Abril de 2020 a las 15:52"

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 469421
$ws.Cells.Item(4, 3).Value = 855
$ws.Cells.Item(4, 4).Value = 25937
$ws.Cells.Item(4, 5).Value = 426774
$ws.Cells.Item(4, 7).Value = 19
$ws.Cells.Item(4, 8).Value = 16710

# Suiza (row 14)
$ws.Cells.Item(14, 2).Value = 24427
$ws.Cells.Item(14, 3).Value = 376
$ws.Cells.Item(14, 5).Value = 12840
$ws.Cells.Item(14, 7).Value = 39
$ws.Cells.Item(14, 8).Value = 987

# Austria (row 19)
$ws.Cells.Item(19, 2).Value = 13492
$ws.Cells.Item(19, 3).Value = 248
$ws.Cells.Item(19, 5).Value = 7109

# Irlanda (row 25)
$ws.Cells.Item(25, 6).Value = 194

# Noruega (row 26)
$ws.Cells.Item(26, 2).Value = 6244
$ws.Cells.Item(26, 3).Value = 25
$ws.Cells.Item(26, 5).Value = 6104

# Singapur jumps ahead of Sudafrica/Bielorrusia/Grecia (rows 52-55)
$ws.Cells.Item(52, 1).Value = "Singapur"
$ws.Cells.Item(52, 2).Value = 2108
$ws.Cells.Item(52, 3).Value = 198
$ws.Cells.Item(52, 4).Value = 492
$ws.Cells.Item(52, 5).Value = 1609
$ws.Cells.Item(52, 6).Value = 29
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = 7

$ws.Cells.Item(53, 1).Value = "Sudafrica"
$ws.Cells.Item(53, 2).Value = 2003
$ws.Cells.Item(53, 3).Value = 69
$ws.Cells.Item(53, 4).Value = 410
$ws.Cells.Item(53, 5).Value = 1569
$ws.Cells.Item(53, 6).Value = 7
$ws.Cells.Item(53, 7).Value = 6
$ws.Cells.Item(53, 8).Value = 24

$ws.Cells.Item(54, 1).Value = "Bielorrusia"
$ws.Cells.Item(54, 2).Value = 1981
$ws.Cells.Item(54, 3).Value = 495
$ws.Cells.Item(54, 4).Value = 169
$ws.Cells.Item(54, 5).Value = 1793
$ws.Cells.Item(54, 6).Value = 72
$ws.Cells.Item(54, 7).Value = 3
$ws.Cells.Item(54, 8).Value = 19

$ws.Cells.Item(55, 1).Value = "Grecia"
$ws.Cells.Item(55, 2).Value = 1955
$ws.Cells.Item(55, 4).Value = 269
$ws.Cells.Item(55, 5).Value = 1599
$ws.Cells.Item(55, 6).Value = 79
$ws.Cells.Item(55, 8).Value = 87

# Islandia jumps ahead of Argelia (rows 58-59)
$ws.Cells.Item(58, 1).Value = "Islandia"
$ws.Cells.Item(58, 2).Value = 1675
$ws.Cells.Item(58, 3).Value = 27
$ws.Cells.Item(58, 4).Value = 751
$ws.Cells.Item(58, 5).Value = 918
$ws.Cells.Item(58, 6).Value = 11
$ws.Cells.Item(58, 8).Value = 6

$ws.Cells.Item(59, 1).Value = "Argelia"
$ws.Cells.Item(59, 2).Value = 1666
$ws.Cells.Item(59, 4).Value = 347
$ws.Cells.Item(59, 5).Value = 1084
$ws.Cells.Item(59, 6).Value = 46
$ws.Cells.Item(59, 8).Value = 235

# Kazajistan (row 76)
$ws.Cells.Item(76, 4).Value = 64
$ws.Cells.Item(76, 5).Value = 729

# Eslovaquia jumps ahead of Crucero/Republica de Macedonia (rows 77-79)
$ws.Cells.Item(77, 1).Value = "Eslovaquia"
$ws.Cells.Item(77, 2).Value = 715
$ws.Cells.Item(77, 3).Value = 14
$ws.Cells.Item(77, 4).Value = 23
$ws.Cells.Item(77, 5).Value = 690
$ws.Cells.Item(77, 6).Value = 5
$ws.Cells.Item(77, 8).Value = 2

$ws.Cells.Item(78, 1).Value = "Crucero"
$ws.Cells.Item(78, 2).Value = 712
$ws.Cells.Item(78, 3).Value = 0
$ws.Cells.Item(78, 4).Value = 619
$ws.Cells.Item(78, 5).Value = 82
$ws.Cells.Item(78, 6).Value = 10
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 11

$ws.Cells.Item(79, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(79, 2).Value = 711
$ws.Cells.Item(79, 3).Value = 48
$ws.Cells.Item(79, 4).Value = 41
$ws.Cells.Item(79, 5).Value = 638
$ws.Cells.Item(79, 6).Value = 15
$ws.Cells.Item(79, 7).Value = 2
$ws.Cells.Item(79, 8).Value = 32

# Ghana (row 99)
$ws.Cells.Item(99, 4).Value = 4
$ws.Cells.Item(99, 5).Value = 368

# San Marino (row 103)
$ws.Cells.Item(103, 2).Value = 344
$ws.Cells.Item(103, 3).Value = 11
$ws.Cells.Item(103, 4).Value = 50
$ws.Cells.Item(103, 5).Value = 260

# Isla de Man (row 117)
$ws.Cells.Item(117, 4).Value = 100
$ws.Cells.Item(117, 5).Value = 89
$ws.Cells.Item(117, 6).Value = 11
